$d = $word.ActiveDocument

# --- Title paragraph: Heading1 -> Title, text "Lesson Plan" -> "Lesson Plan: "
$titlePara = $d.Paragraphs(1)
$titlePara.Style = "Title"
$titlePara.Range.Text = "Lesson Plan: "

# --- Insert the lesson-details table right after the title paragraph
$insertionRange = $d.Paragraphs(1).Range
$insertionRange.Collapse(0)
$tbl = $d.Tables.Add($insertionRange, 4, 2)

# Column widths: 5396 twips = 269.8 points each (content width once right
# margin becomes 724, matching the page's usable width of 10792 twips)
$tbl.Columns(1).Width = 269.8
$tbl.Columns(2).Width = 269.8

$rows = @(
    @("Subject", "Physical Education"),
    @("Date", ""),
    @("Grade/Level", "8"),
    @("Suggested Lesson Time", "30 minutes")
)

for ($i = 1; $i -le 4; $i++) {
    $labelCell = $tbl.Rows($i).Cells(1)
    $valueCell = $tbl.Rows($i).Cells(2)

    $labelCell.Range.Paragraphs(1).Style = "Normal"
    $labelCell.Range.Text = $rows[$i - 1][0]

    $valueCell.Range.Paragraphs(1).Style = "Normal"
    $valueText = $rows[$i - 1][1]
    if ($valueText -ne "") {
        $valueCell.Range.Text = $valueText
    }
}

# --- Section margins: right margin 180 -> 724 twips (9 -> 36.2pt)
$sec = $d.Sections(1)
$sec.PageSetup.RightMargin = 36.2
